# Apply revised statistics to the "Panel B (E-mini Futures) - Avg Daily Volume"
# block (row 26), its "Diff_Vol (Ann - Day)" row (27) and "# Obs" row (28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: Avg Daily Volume
$ws.Range("D26").Value = 467865.686746988
$ws.Range("E26").Value = 678248.848442145
$ws.Range("G26").Value = 5323
$ws.Range("H26").Value = 727485.5
$ws.Range("I26").Value = 83
$ws.Range("J26").Value = 661285.2530120482
$ws.Range("K26").Value = 861433.7148268499
$ws.Range("M26").Value = 519761
$ws.Range("N26").Value = 1044893
$ws.Range("O26").Value = 83
$ws.Range("P26").Value = 627804.7951807228
$ws.Range("Q26").Value = 774264.2255404798
$ws.Range("S26").Value = 553930
$ws.Range("T26").Value = 1030945.5
$ws.Range("U26").Value = 83
$ws.Range("V26").Value = 562209.4939759036
$ws.Range("W26").Value = 780689.979638749
$ws.Range("Y26").Value = 15134
$ws.Range("Z26").Value = 941555
$ws.Range("AA26").Value = 83
$ws.Range("AB26").Value = 535375.1445783132
$ws.Range("AC26").Value = 703681.3614342463
$ws.Range("AE26").Value = 169572
$ws.Range("AF26").Value = 838115
$ws.Range("AG26").Value = 83

# Row 27: Diff_Vol (Ann - Day)
$ws.Range("D27").Value = 159939.1084337349
$ws.Range("J27").Value = -33480.4578313253
$ws.Range("V27").Value = 65595.30120481928
$ws.Range("AB27").Value = 92429.65060240965

# Row 28: # Obs
$ws.Range("D28").Value = 83
$ws.Range("J28").Value = 83
$ws.Range("P28").Value = 83
$ws.Range("V28").Value = 83
$ws.Range("AB28").Value = 83
